$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("I2").Value = 3933
$ws.Range("I3").Value = 4076
$ws.Range("I4").Value = 950
$ws.Range("I5").Value = 376
$ws.Range("I6").Value = 4537
$ws.Range("I7").Value = 13872

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("I2").Value = 38
$ws.Range("I7").Value = 155

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("I2").Value = 42
$ws.Range("I7").Value = 150

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("I4").Value = 3
$ws.Range("I7").Value = 76

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("I3").Value = 144
$ws.Range("I6").Value = 121
$ws.Range("I7").Value = 446

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("I6").Value = 70
$ws.Range("I7").Value = 258

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("I3").Value = 192
$ws.Range("I6").Value = 178
$ws.Range("I7").Value = 538

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("I3").Value = 36
$ws.Range("I7").Value = 136

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("I2").Value = 43
$ws.Range("I7").Value = 117

$ws = $wb.Worksheets.Item('New City')
$ws.Range("I2").Value = 101
$ws.Range("I4").Value = 14
$ws.Range("I7").Value = 310

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("I2").Value = 120
$ws.Range("I5").Value = 44
$ws.Range("I7").Value = 437
$ws.Range("I8").Value = 834
$ws.Range("I10").Value = 93
$ws.Range("I11").Value = 212
$ws.Range("I14").Value = 76
$ws.Range("I15").Value = 166
$ws.Range("I18").Value = 95
$ws.Range("I19").Value = 382
$ws.Range("I20").Value = 341
$ws.Range("I29").Value = 898
$ws.Range("I31").Value = 136
$ws.Range("I33").Value = 634
$ws.Range("I37").Value = 446
$ws.Range("I42").Value = 481
$ws.Range("I43").Value = 120
$ws.Range("I48").Value = 184
$ws.Range("I51").Value = 135
$ws.Range("I54").Value = 317
$ws.Range("I60").Value = 68
$ws.Range("I63").Value = 52
$ws.Range("I65").Value = 310
$ws.Range("I67").Value = 538
$ws.Range("I71").Value = 42
$ws.Range("I73").Value = 115
$ws.Range("I76").Value = 210
$ws.Range("I77").Value = 76
$ws.Range("I79").Value = 374
$ws.Range("I83").Value = 279
$ws.Range("I84").Value = 117
$ws.Range("I85").Value = 626
$ws.Range("I88").Value = 124
$ws.Range("I89").Value = 155
$ws.Range("I91").Value = 173
$ws.Range("I95").Value = 222
$ws.Range("I96").Value = 150
$ws.Range("I98").Value = 92
$ws.Range("I99").Value = 258
$ws.Range("I101").Value = 13872

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("I3").Value = 107
$ws.Range("I7").Value = 279

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("I2").Value = 81
$ws.Range("I7").Value = 222

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("I2").Value = 149
$ws.Range("I3").Value = 234
$ws.Range("I6").Value = 197
$ws.Range("I7").Value = 634

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("I2").Value = 72
$ws.Range("I3").Value = 63
$ws.Range("I6").Value = 159
$ws.Range("I7").Value = 317

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("I2").Value = 265
$ws.Range("I3").Value = 307
$ws.Range("I5").Value = 37
$ws.Range("I6").Value = 245
$ws.Range("I7").Value = 898

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("I2").Value = 142
$ws.Range("I3").Value = 110
$ws.Range("I7").Value = 382

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("I2").Value = 22
$ws.Range("I6").Value = 108
$ws.Range("I7").Value = 184

$ws = $wb.Worksheets.Item('River North')
$ws.Range("I6").Value = 87
$ws.Range("I7").Value = 210

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("I2").Value = 162
$ws.Range("I3").Value = 252
$ws.Range("I6").Value = 155
$ws.Range("I7").Value = 626

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("I5").Value = 17
$ws.Range("I6").Value = 131
$ws.Range("I7").Value = 481

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("I2").Value = 32
$ws.Range("I7").Value = 93

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("I3").Value = 61
$ws.Range("I7").Value = 173

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("I4").Value = 24
$ws.Range("I6").Value = 112
$ws.Range("I7").Value = 374

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("I3").Value = 104
$ws.Range("I6").Value = 108
$ws.Range("I7").Value = 341

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("I6").Value = 40
$ws.Range("I7").Value = 95

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("I2").Value = 53
$ws.Range("I6").Value = 59
$ws.Range("I7").Value = 166

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("I6").Value = 59
$ws.Range("I7").Value = 92

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("I6").Value = 51
$ws.Range("I7").Value = 212

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("I2").Value = 41
$ws.Range("I3").Value = 32
$ws.Range("I7").Value = 115

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("I3").Value = 42
$ws.Range("I7").Value = 120

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("I3").Value = 45
$ws.Range("I7").Value = 124

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("I3").Value = 234
$ws.Range("I6").Value = 270
$ws.Range("I7").Value = 834

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range("I3").Value = 14
$ws.Range("I6").Value = 21
$ws.Range("I7").Value = 44

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("I4").Value = 15
$ws.Range("I7").Value = 135

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("I2").Value = 18
$ws.Range("I7").Value = 68

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("I6").Value = 69
$ws.Range("I7").Value = 120

$ws = $wb.Worksheets.Item('Oakland')
$ws.Range("I4").Value = 3
$ws.Range("I7").Value = 42

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("I2").Value = 23
$ws.Range("I7").Value = 76

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("I2").Value = 150
$ws.Range("I3").Value = 134
$ws.Range("I7").Value = 437

$ws = $wb.Worksheets.Item('Beverly')
$ws.Range("I2").Value = 6
$ws.Range("I6").Value = 12
